$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng1 = $ws.Range("G40")
$rng1.Interior.PatternThemeColor = 5
Write-Host "done check patternthemecolor alone"
